# Updates the cryptos list values (prices and 1h volume changes) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.633.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.397.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.13%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'254.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.63%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'652.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.12%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.81%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.431"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.25%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +6.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D11").Value = "'3.392.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.212"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'41.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.33%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +16.65%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0000260"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.46%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'97.215.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.34%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.024.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'8.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +27.89%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.404.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.50%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.42%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.509"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +49.25%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +10.79%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.98%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'508.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.14%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0000206"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'NEARProtocol"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'6.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.51%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Litecoin"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'99.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +10.42%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'12.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.31%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.576.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.21%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.154"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.74%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.207"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'11.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.21%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.996"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.14%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.567"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +15.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'29.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.71%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +13.89%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'7.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.99%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +11.91%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'523.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.61%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.153"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.17%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.12%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.856"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +5.72%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0423"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +17.68%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.18%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'3.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.84%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'5.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +10.95%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Cosmos"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'8.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +10.14%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'USDe"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.08%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +8.05%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.93%  "
$ws.Range("E51").Style = "Normal"
